$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Formula = "=0.1/2"
$ws.Range("G7").Formula = "=0.3/2"
$ws.Range("H7").Formula = "=(E7+G7)/2"
$ws.Range("Q7").Formula = "=(0.38 + 0.75)/2"
$ws.Range("T7").Formula = "=(Q7+S7)/2"

$ws.Range("E8").Formula = "=0.45/2"
$ws.Range("G8").Formula = "=(0.9+0.45)/2"
$ws.Range("H8").Formula = "=(E8+G8)/2"
$ws.Range("R8").Value = "plannedaardvarkdrinksbloodyriverbed"
$ws.Range("S8").Value = 1
$ws.Range("T8").Value = 1

$ws.Range("C8").Select()
$ws.Application.ActiveWindow.ScrollColumn = 3
